$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.379.12'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.77%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.883.14'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.697'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '246.53'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.43'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.70%  '
$ws.Range('E9').Value = '  +2.77%  '
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '13.51'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.159.21'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.773'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +8.29%  '
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.870.63'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.380.98'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '73.49'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0827'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '244.81'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.98%  '
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.18'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.61'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +8.45%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.17'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.56'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.66'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.32'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.128'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.76%  '
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('E31').Value = '  +0.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.89'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.19'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  -12.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.852'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.09%  '
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('E38').Value = '  +11.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.30'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('E40').Value = '  +3.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.41'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.07'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('E43').Value = '  +2.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.309.32'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0811'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.16%  '
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.12'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('E49').Value = '  -2.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.32'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.063.36'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.36%  '
